# Update NATMI ligand-receptor (Col4a2-Cd93) output with refreshed TPM-based values.
# Ligand/receptor average & total expression, their derived specificities, and the
# resulting edge weights/specificities all change for every sending x target cluster
# combination in rows 2-17 (columns G,H,I,J,M,N,O,P,Q,R,S,T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 231.1867016666667
$ws.Range("H2").Value2 = 693.560105
$ws.Range("I2").Value2 = 0.5918542142967572
$ws.Range("J2").Value2 = 0.5918542142967572
$ws.Range("M2").Value2 = 209.26237
$ws.Range("N2").Value2 = 627.78711
$ws.Range("O2").Value2 = 0.8127157202241573
$ws.Range("P2").Value2 = 0.8127157202241573
$ws.Range("Q2").Value2 = 48378.67710324962
$ws.Range("R2").Value2 = 435408.0939292465
$ws.Range("S2").Value2 = 0.4810092240398918
$ws.Range("T2").Value2 = 0.4810092240398918

$ws.Range("G3").Value2 = 231.1867016666667
$ws.Range("H3").Value2 = 693.560105
$ws.Range("I3").Value2 = 0.5918542142967572
$ws.Range("J3").Value2 = 0.5918542142967572
$ws.Range("M3").Value2 = 0.9848756666666668
$ws.Range("N3").Value2 = 2.954627
$ws.Range("O3").Value2 = 0.003824977881910862
$ws.Range("P3").Value2 = 0.003824977881910862
$ws.Range("Q3").Value2 = 227.6901569284261
$ws.Range("R3").Value2 = 2049.211412355835
$ws.Range("S3").Value2 = 0.002263829279000828
$ws.Range("T3").Value2 = 0.002263829279000828

$ws.Range("G4").Value2 = 231.1867016666667
$ws.Range("H4").Value2 = 693.560105
$ws.Range("I4").Value2 = 0.5918542142967572
$ws.Range("J4").Value2 = 0.5918542142967572
$ws.Range("M4").Value2 = 1.763846666666667
$ws.Range("N4").Value2 = 5.291539999999999
$ws.Range("O4").Value2 = 0.006850280411451801
$ws.Range("P4").Value2 = 0.006850280411451801
$ws.Range("Q4").Value2 = 407.7778931124111
$ws.Range("R4").Value2 = 3670.0010380117
$ws.Range("S4").Value2 = 0.004054367330632272
$ws.Range("T4").Value2 = 0.004054367330632272

$ws.Range("G5").Value2 = 231.1867016666667
$ws.Range("H5").Value2 = 693.560105
$ws.Range("I5").Value2 = 0.5918542142967572
$ws.Range("J5").Value2 = 0.5918542142967572
$ws.Range("M5").Value2 = 45.474231
$ws.Range("N5").Value2 = 136.422693
$ws.Range("O5").Value2 = 0.1766090214824801
$ws.Range("P5").Value2 = 0.1766090214824801
$ws.Range("Q5").Value2 = 10513.03747571808
$ws.Range("R5").Value2 = 94617.33728146275
$ws.Range("S5").Value2 = 0.1045267936472324
$ws.Range("T5").Value2 = 0.1045267936472324

$ws.Range("I6").Value2 = 0.2653686241974456
$ws.Range("J6").Value2 = 0.2653686241974456
$ws.Range("M6").Value2 = 209.26237
$ws.Range("N6").Value2 = 627.78711
$ws.Range("O6").Value2 = 0.8127157202241573
$ws.Range("P6").Value2 = 0.8127157202241573
$ws.Range("Q6").Value2 = 21691.4616357614
$ws.Range("R6").Value2 = 195223.1547218526
$ws.Range("S6").Value2 = 0.2156692525395207
$ws.Range("T6").Value2 = 0.2156692525395207

$ws.Range("I7").Value2 = 0.2653686241974456
$ws.Range("J7").Value2 = 0.2653686241974456
$ws.Range("M7").Value2 = 0.9848756666666668
$ws.Range("N7").Value2 = 2.954627
$ws.Range("O7").Value2 = 0.003824977881910862
$ws.Range("P7").Value2 = 0.003824977881910862
$ws.Range("R7").Value2 = 918.8012859428782
$ws.Range("S7").Value2 = 0.001015029118108345
$ws.Range("T7").Value2 = 0.001015029118108345

$ws.Range("I8").Value2 = 0.2653686241974456
$ws.Range("J8").Value2 = 0.2653686241974456
$ws.Range("M8").Value2 = 1.763846666666667
$ws.Range("N8").Value2 = 5.291539999999999
$ws.Range("O8").Value2 = 0.006850280411451801
$ws.Range("P8").Value2 = 0.006850280411451801
$ws.Range("Q8").Value2 = 182.8346505937289
$ws.Range("R8").Value2 = 1645.51185534356
$ws.Range("S8").Value2 = 0.001817849488153676
$ws.Range("T8").Value2 = 0.001817849488153676

$ws.Range("I9").Value2 = 0.2653686241974456
$ws.Range("J9").Value2 = 0.2653686241974456
$ws.Range("M9").Value2 = 45.474231
$ws.Range("N9").Value2 = 136.422693
$ws.Range("O9").Value2 = 0.1766090214824801
$ws.Range("P9").Value2 = 0.1766090214824801
$ws.Range("Q9").Value2 = 4713.711964326178
$ws.Range("R9").Value2 = 42423.4076789356
$ws.Range("S9").Value2 = 0.04686649305166285
$ws.Range("T9").Value2 = 0.04686649305166285

$ws.Range("G10").Value2 = 54.625754
$ws.Range("H10").Value2 = 163.877262
$ws.Range("I10").Value2 = 0.1398457717548702
$ws.Range("J10").Value2 = 0.1398457717548702
$ws.Range("M10").Value2 = 209.26237
$ws.Range("N10").Value2 = 627.78711
$ws.Range("O10").Value2 = 0.8127157202241573
$ws.Range("P10").Value2 = 0.8127157202241573
$ws.Range("Q10").Value2 = 11431.11474507698
$ws.Range("R10").Value2 = 102880.0327056928
$ws.Range("S10").Value2 = 0.1136548571120625
$ws.Range("T10").Value2 = 0.1136548571120624

$ws.Range("G11").Value2 = 54.625754
$ws.Range("H11").Value2 = 163.877262
$ws.Range("I11").Value2 = 0.1398457717548702
$ws.Range("J11").Value2 = 0.1398457717548702
$ws.Range("M11").Value2 = 0.9848756666666668
$ws.Range("N11").Value2 = 2.954627
$ws.Range("O11").Value2 = 0.003824977881910862
$ws.Range("P11").Value2 = 0.003824977881910862
$ws.Range("Q11").Value2 = 53.79957588791934
$ws.Range("R11").Value2 = 484.196182991274
$ws.Range("S11").Value2 = 0.0005349069838411334
$ws.Range("T11").Value2 = 0.0005349069838411333

$ws.Range("G12").Value2 = 54.625754
$ws.Range("H12").Value2 = 163.877262
$ws.Range("I12").Value2 = 0.1398457717548702
$ws.Range("J12").Value2 = 0.1398457717548702
$ws.Range("M12").Value2 = 1.763846666666667
$ws.Range("N12").Value2 = 5.291539999999999
$ws.Range("O12").Value2 = 0.006850280411451801
$ws.Range("P12").Value2 = 0.006850280411451801
$ws.Range("Q12").Value2 = 96.35145410705333
$ws.Range("R12").Value2 = 867.1630869634799
$ws.Range("S12").Value2 = 0.0009579827508767471
$ws.Range("T12").Value2 = 0.0009579827508767469

$ws.Range("G13").Value2 = 54.625754
$ws.Range("H13").Value2 = 163.877262
$ws.Range("I13").Value2 = 0.1398457717548702
$ws.Range("J13").Value2 = 0.1398457717548702
$ws.Range("M13").Value2 = 45.474231
$ws.Range("N13").Value2 = 136.422693
$ws.Range("O13").Value2 = 0.1766090214824801
$ws.Range("P13").Value2 = 0.1766090214824801
$ws.Range("Q13").Value2 = 2484.064155945174
$ws.Range("R13").Value2 = 22356.57740350656
$ws.Range("S13").Value2 = 0.02469802490808988
$ws.Range("T13").Value2 = 0.02469802490808988

$ws.Range("G14").Value2 = 1.145042666666667
$ws.Range("H14").Value2 = 3.435128
$ws.Range("I14").Value2 = 0.002931389750926909
$ws.Range("J14").Value2 = 0.002931389750926909
$ws.Range("M14").Value2 = 209.26237
$ws.Range("N14").Value2 = 627.78711
$ws.Range("O14").Value2 = 0.8127157202241573
$ws.Range("P14").Value2 = 0.8127157202241573
$ws.Range("Q14").Value2 = 239.6143421777867
$ws.Range("R14").Value2 = 2156.52907960008
$ws.Range("S14").Value2 = 0.002382386532682276
$ws.Range("T14").Value2 = 0.002382386532682276

$ws.Range("G15").Value2 = 1.145042666666667
$ws.Range("H15").Value2 = 3.435128
$ws.Range("I15").Value2 = 0.002931389750926909
$ws.Range("J15").Value2 = 0.002931389750926909
$ws.Range("M15").Value2 = 0.9848756666666668
$ws.Range("N15").Value2 = 2.954627
$ws.Range("O15").Value2 = 0.003824977881910862
$ws.Range("P15").Value2 = 0.003824977881910862
$ws.Range("Q15").Value2 = 1.127724659695111
$ws.Range("R15").Value2 = 10.149521937256
$ws.Range("S15").Value2 = 0.00001121250096055562
$ws.Range("T15").Value2 = 0.00001121250096055562

$ws.Range("G16").Value2 = 1.145042666666667
$ws.Range("H16").Value2 = 3.435128
$ws.Range("I16").Value2 = 0.002931389750926909
$ws.Range("J16").Value2 = 0.002931389750926909
$ws.Range("M16").Value2 = 1.763846666666667
$ws.Range("N16").Value2 = 5.291539999999999
$ws.Range("O16").Value2 = 0.006850280411451801
$ws.Range("P16").Value2 = 0.006850280411451801
$ws.Range("Q16").Value2 = 2.019679690791111
$ws.Range("R16").Value2 = 18.17711721712
$ws.Range("S16").Value2 = 0.00002008084178910518
$ws.Range("T16").Value2 = 0.00002008084178910518

$ws.Range("G17").Value2 = 1.145042666666667
$ws.Range("H17").Value2 = 3.435128
$ws.Range("I17").Value2 = 0.002931389750926909
$ws.Range("J17").Value2 = 0.002931389750926909
$ws.Range("M17").Value2 = 45.474231
$ws.Range("N17").Value2 = 136.422693
$ws.Range("O17").Value2 = 0.1766090214824801
$ws.Range("P17").Value2 = 0.1766090214824801
$ws.Range("Q17").Value2 = 52.069934728856
$ws.Range("R17").Value2 = 468.6294125597039
$ws.Range("S17").Value2 = 0.0005177098754949724
$ws.Range("T17").Value2 = 0.0005177098754949724
